$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 0.9999932936434416
$ws.Range("E2").Value = 0.9999932936434416

# Row 3
$ws.Range("D3").Value = 0.9999070183706054
$ws.Range("E3").Value = 0.9999070183706054

# Row 4
$ws.Range("D4").Value = 0.9398987122356667
$ws.Range("E4").Value = 0.9398987122356667

# Row 5
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = 0.07176928247659721
$ws.Range("E5").Value = 0.07176928247659721

# Row 6
$ws.Range("D6").Value = 0.3105843576357338
$ws.Range("E6").Value = 0.3105843576357338

# Row 7
$ws.Range("D7").Value = 0.9966264512541333
$ws.Range("E7").Value = 0.003373548745866706

# Row 8
$ws.Range("D8").Value = 0.9999999999999998
$ws.Range("E8").Value = [double]"2.220446049250313E-16"

# Row 9
$ws.Range("D9").Value = 0.03595781548410305
$ws.Range("E9").Value = 0.964042184515897

# Row 10
$ws.Range("D10").Value = 0.9999999999894822
$ws.Range("E10").Value = [double]"1.051780884608888E-11"

# Row 11
$ws.Range("D11").Value = 0.08974363309369945
$ws.Range("E11").Value = 0.9102563669063005
$ws.Range("F11").Value = 3.019325733184814
$ws.Range("G11").Value = 0.5
